# Generate Report for Archive
#
# The shared "handoff" status text is updated from "Ready for handoff" to
# "In Translation" everywhere it appears (Overview!E2:F2, zh-cn!C2, de-de!C2).
# Updating every cell that references the old text lets the workbook's string
# table collapse back down to a single shared entry for the new text, which
# keeps cell E2/F2 (Overview) and C2 (zh-cn / de-de) all pointing at the same
# string - exactly mirroring the source diff's in-place shared-string edit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# The status column(s) narrow to fit the shorter replacement text
# (target stored width ~= 13.41 characters). The host snaps ColumnWidth
# onto its internal pixel grid, so 12.5 is the input that lands on the
# grid point closest to the target width.
$newStatusColWidth = 12.5
$overview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$zhcn.Columns.Item(3).ColumnWidth = $newStatusColWidth
$dede.Columns.Item(3).ColumnWidth = $newStatusColWidth
